$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that gets bumped by one day
# (45181 -> 45182) for every data row (rows 2 through 236).
$ws.Range("C2:C236").Value = 45182
